# Apply "Ran code for averaged intensities on spiral schemes":
# - Relabel/renumber rows so the table keeps HKL-index continuity while
#   inserting the 3 new Spiral-* schemes right after Gaussian-Quadrature,
#   and append the 3 rows (HexGrid family) that got pushed off the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "ND Single"
$ws.Range("C3").Value = 1.05
$ws.Range("D3").Value = 0.6899999999999999
$ws.Range("E3").Value = 1.07
$ws.Range("F3").Value = 1.05
$ws.Range("G3").Value = 0.86
$ws.Range("H3").Value = 1.2
$ws.Range("I3").Value = 1.07
$ws.Range("J3").Value = 0.6899999999999999
$ws.Range("K3").Value = 0.88
$ws.Range("L3").Value = 0.9650000000000001
$ws.Range("M3").Value = 0.9900000000000001

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "RD Single"
$ws.Range("C4").Value = 1.02
$ws.Range("D4").Value = 0.89
$ws.Range("E4").Value = 1.01
$ws.Range("F4").Value = 1.02
$ws.Range("G4").Value = 0.93
$ws.Range("H4").Value = 1.08
$ws.Range("I4").Value = 1.02
$ws.Range("J4").Value = 0.89
$ws.Range("K4").Value = 0.95
$ws.Range("L4").Value = 0.985
$ws.Range("M4").Value = 0.9916666666666666

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "TD Single"
$ws.Range("C5").Value = 1.03
$ws.Range("D5").Value = 0.88
$ws.Range("E5").Value = 1.02
$ws.Range("F5").Value = 1.03
$ws.Range("G5").Value = 0.9399999999999999
$ws.Range("H5").Value = 1.05
$ws.Range("I5").Value = 1.02
$ws.Range("J5").Value = 0.88
$ws.Range("K5").Value = 0.95
$ws.Range("L5").Value = 0.99
$ws.Range("M5").Value = 0.9899999999999999

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Morris"
$ws.Range("C6").Value = 0.98
$ws.Range("D6").Value = 1.09
$ws.Range("E6").Value = 0.98
$ws.Range("F6").Value = 0.98
$ws.Range("G6").Value = 1.04
$ws.Range("H6").Value = 0.9399999999999999
$ws.Range("I6").Value = 0.98
$ws.Range("J6").Value = 1.09
$ws.Range("K6").Value = 1.035
$ws.Range("L6").Value = 1.0075
$ws.Range("M6").Value = 1.001666666666667

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Ring Perpendicular to ND"
$ws.Range("C7").Value = 1.028767123287671
$ws.Range("D7").Value = 0.8827397260273973
$ws.Range("E7").Value = 1.015753424657534
$ws.Range("F7").Value = 1.028767123287671
$ws.Range("G7").Value = 0.9363013698630137
$ws.Range("H7").Value = 1.069315068493151
$ws.Range("I7").Value = 1.02
$ws.Range("J7").Value = 0.8827397260273973
$ws.Range("K7").Value = 0.9492465753424657
$ws.Range("L7").Value = 0.9890068493150683
$ws.Range("M7").Value = 0.9921461187214611

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Ring Perpendicular to RD"
$ws.Range("C8").Value = 1.008947368421053
$ws.Range("D8").Value = 0.9457894736842105
$ws.Range("E8").Value = 1.008421052631579
$ws.Range("F8").Value = 1.008947368421053
$ws.Range("G8").Value = 0.9694736842105263
$ws.Range("H8").Value = 1.026842105263158
$ws.Range("I8").Value = 1.007894736842105
$ws.Range("J8").Value = 0.9457894736842105
$ws.Range("K8").Value = 0.9771052631578947
$ws.Range("L8").Value = 0.9930263157894736
$ws.Range("M8").Value = 0.9945614035087718

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Ring Perpendicular to TD"
$ws.Range("C9").Value = 1.004210526315789
$ws.Range("D9").Value = 0.9473684210526315
$ws.Range("E9").Value = 1.009473684210526
$ws.Range("F9").Value = 1.004210526315789
$ws.Range("G9").Value = 0.9710526315789474
$ws.Range("H9").Value = 1.034210526315789
$ws.Range("I9").Value = 1.007368421052631
$ws.Range("J9").Value = 0.9473684210526315
$ws.Range("K9").Value = 0.9784210526315789
$ws.Range("L9").Value = 0.9913157894736842
$ws.Range("M9").Value = 0.9956140350877193

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.013849313885671
$ws.Range("D10").Value = 0.9229541462827703
$ws.Range("E10").Value = 1.011873886281002
$ws.Range("F10").Value = 1.013849313885671
$ws.Range("G10").Value = 0.960652774054112
$ws.Range("H10").Value = 1.041031941877974
$ws.Range("I10").Value = 1.012512863396901
$ws.Range("J10").Value = 0.9229541462827703
$ws.Range("K10").Value = 0.9674140162818861
$ws.Range("L10").Value = 0.9906316650837783
$ws.Range("M10").Value = 0.9938124876297384

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 1.00014273700017
$ws.Range("D11").Value = 0.9684239342159544
$ws.Range("E11").Value = 1.004894700758721
$ws.Range("F11").Value = 1.00014273700017
$ws.Range("G11").Value = 0.9806277731655929
$ws.Range("H11").Value = 1.01989118052493
$ws.Range("I11").Value = 1.002346189866539
$ws.Range("J11").Value = 0.9684239342159544
$ws.Range("K11").Value = 0.9866593174873377
$ws.Range("L11").Value = 0.9934010272437537
$ws.Range("M11").Value = 0.9960544192553177
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 1.00011158192021
$ws.Range("D12").Value = 0.9686591048010311
$ws.Range("E12").Value = 1.004870083695703
$ws.Range("F12").Value = 1.00011158192021
$ws.Range("G12").Value = 0.9807042768308752
$ws.Range("H12").Value = 1.01976066362106
$ws.Range("I12").Value = 1.00228706557193
$ws.Range("J12").Value = 0.9686591048010311
$ws.Range("K12").Value = 0.9867645942483673
$ws.Range("L12").Value = 0.9934380880842888
$ws.Range("M12").Value = 0.9960654627401352
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 1.000159092737074
$ws.Range("D13").Value = 0.9684999677098959
$ws.Range("E13").Value = 1.004879091491546
$ws.Range("F13").Value = 1.000159092737074
$ws.Range("G13").Value = 0.9805882720709517
$ws.Range("H13").Value = 1.01989030710423
$ws.Range("I13").Value = 1.002322878855956
$ws.Range("J13").Value = 0.9684999677098959
$ws.Range("K13").Value = 0.986689529600721
$ws.Range("L13").Value = 0.9934243111688974
$ws.Range("M13").Value = 0.9960566016616088
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 1.028328
$ws.Range("D14").Value = 0.8189080000000011
$ws.Range("E14").Value = 1.039387999999999
$ws.Range("F14").Value = 1.028328
$ws.Range("G14").Value = 0.9158160000000002
$ws.Range("H14").Value = 1.115243999999999
$ws.Range("I14").Value = 1.038343999999999
$ws.Range("J14").Value = 0.8189080000000011
$ws.Range("K14").Value = 0.9291480000000002
$ws.Range("L14").Value = 0.9787380000000001
$ws.Range("M14").Value = 0.9926713333333331

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1.05
$ws.Range("D15").Value = 0.6899999999999999
$ws.Range("E15").Value = 1.07
$ws.Range("F15").Value = 1.05
$ws.Range("G15").Value = 0.86
$ws.Range("H15").Value = 1.2
$ws.Range("I15").Value = 1.07
$ws.Range("J15").Value = 0.6899999999999999
$ws.Range("K15").Value = 0.88
$ws.Range("L15").Value = 0.9650000000000001
$ws.Range("M15").Value = 0.9900000000000001

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 1.028069541068803
$ws.Range("D16").Value = 0.819464077107205
$ws.Range("E16").Value = 1.039250154496
$ws.Range("F16").Value = 1.028069541068803
$ws.Range("G16").Value = 0.9160449341440041
$ws.Range("H16").Value = 1.114363213004794
$ws.Range("I16").Value = 1.038279307468797
$ws.Range("J16").Value = 0.819464077107205
$ws.Range("K16").Value = 0.9293571158016023
$ws.Range("L16").Value = 0.9787133284352026
$ws.Range("M16").Value = 0.9925785378816006

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.995287826776786
$ws.Range("D17").Value = 0.9962186513699979
$ws.Range("E17").Value = 0.9969416023052503
$ws.Range("F17").Value = 0.995287826776786
$ws.Range("G17").Value = 0.9957110568673143
$ws.Range("H17").Value = 0.9971101070422019
$ws.Range("I17").Value = 0.9960247792402397
$ws.Range("J17").Value = 0.9962186513699979
$ws.Range("K17").Value = 0.9965801268376241
$ws.Range("L17").Value = 0.995933976807205
$ws.Range("M17").Value = 0.9962156706002984
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9952648451015624
$ws.Range("D18").Value = 1.003315778163301
$ws.Range("E18").Value = 0.9945721139656404
$ws.Range("F18").Value = 0.9952648451015624
$ws.Range("G18").Value = 1.001215824023769
$ws.Range("H18").Value = 0.9912214706489111
$ws.Range("I18").Value = 0.9952533546812806
$ws.Range("J18").Value = 1.003315778163301
$ws.Range("K18").Value = 0.9989439460644705
$ws.Range("L18").Value = 0.9971043955830166
$ws.Range("M18").Value = 0.9968072310974107
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.990152256744693
$ws.Range("D19").Value = 1.025675533047859
$ws.Range("E19").Value = 0.9895289401362373
$ws.Range("F19").Value = 0.990152256744693
$ws.Range("G19").Value = 1.012469911095
$ws.Range("H19").Value = 0.9742612169895541
$ws.Range("I19").Value = 0.9883970631097814
$ws.Range("J19").Value = 1.025675533047859
$ws.Range("K19").Value = 1.007602236592048
$ws.Range("L19").Value = 0.9988772466683706
$ws.Range("M19").Value = 0.9967474868538542
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
